$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.608.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.059.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.12%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.050.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.47%  "

$ws.Range("E9").Value = "  +1.72%  "

$ws.Range("E10").Value = "  +5.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.10%  "

$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("E13").Value = "  +4.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.58%  "

$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.559.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.051.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "61.556.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "450.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +4.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("E33").Value = "  +6.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0806"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.36%  "

$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "412.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0363"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.784.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.94%  "

$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "37.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.69%  "

$ws.Range("E47").Value = "  +5.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("E50").Value = "  +1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.77%  "
